$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / title cell
$ws.Range("A1").Value = "cdwa1"

# Data table: columns B..G, rows 4..23 (run 0..19), plus row24 (Mean, columns B..F)
$data = @(
    @(140,   -630,  -960,  -800,  770,  "exp3_cdwa1_w3_C4_r0.bag"),
    @(-900, -1280, -1500,  1490, -430,  "exp3_cdwa1_w3_C4_r1.bag"),
    @(-380,     0,   630,   790,    0,  "exp3_cdwa1_w3_C4_r2.bag"),
    @(700,   1320,  1490,   720, 1490,  "exp3_cdwa1_w3_C4_r3.bag"),
    @(-390,   -20,   170,   330,    0,  "exp3_cdwa1_w3_C4_r4.bag"),
    @(-540,  -360, -1500,  1490,  -50,  "exp3_cdwa1_w3_C8_r0.bag"),
    @(0,     -250,  -230,  -660,  100,  "exp3_cdwa1_w3_C8_r1.bag"),
    @(-230,  -250,     0,   340,  140,  "exp3_cdwa1_w3_C8_r2.bag"),
    @(-620,  -490, -1140, -1500, -140,  "exp3_cdwa1_w3_C8_r3.bag"),
    @(-560,  -530,  -440,   250,  -30,  "exp3_cdwa1_w3_C8_r4.bag"),
    @(-680,  -840, -1500,     0, -200,  "exp3_cdwa1_w4_C4_r0.bag"),
    @(1490,  1490,  1490,  1360,    0,  "exp3_cdwa1_w4_C4_r1.bag"),
    @(-220,     0,   940,  1290,   80,  "exp3_cdwa1_w4_C4_r2.bag"),
    @(-1500, -1500, 1490,  1490, -1500, "exp3_cdwa1_w4_C4_r3.bag"),
    @(-420,    30,   610,   890,    0,  "exp3_cdwa1_w4_C4_r4.bag"),
    @(-640,  -660,  -600,     0, -120,  "exp3_cdwa1_w4_C8_r0.bag"),
    @(1490,  1490,   580,   900,    0,  "exp3_cdwa1_w4_C8_r1.bag"),
    @(-410,   -20,   700,  1120,    0,  "exp3_cdwa1_w4_C8_r2.bag"),
    @(-600,  -450,  -160,   730,  -30,  "exp3_cdwa1_w4_C8_r3.bag"),
    @(-460,  -300,     0,     0,    0,  "exp3_cdwa1_w4_C8_r4.bag")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

# Mean row (row 24), columns B..F only
$ws.Cells.Item(24, 2).Value = -236
$ws.Cells.Item(24, 3).Value = -162
$ws.Cells.Item(24, 4).Value = 3
$ws.Cells.Item(24, 5).Value = 511
$ws.Cells.Item(24, 6).Value = 4

$wb.Save()
